# Trading update: 2026-02-17 20:17:42
# Appends the newest MarketMaking trade (Trade #35, opened 20:17:26) to both
# the "All Trades" log (next blank row) and the "MarketMaking" strategy
# worksheet (next blank row).

$wb = $excel.ActiveWorkbook

# ---- "All Trades" sheet: currently A1:Q35 -> new trade lands on row 36 ----
$ws = $wb.Worksheets.Item("All Trades")

$ws.Cells.Item(36, 1).Value = 35

# Date string "2026-02-17" -- force text so the engine does not smart-parse
# it into a date serial, then drop the temporary style so the cell stays
# unstyled like the rest of the sheet.
$ws.Cells.Item(36, 2).NumberFormat = "@"
$ws.Cells.Item(36, 2).Value = "2026-02-17"
$ws.Cells.Item(36, 2).ClearFormats()

$ws.Cells.Item(36, 3).Value = "20:17:26"
$ws.Cells.Item(36, 4).Value = "MarketMaking"
$ws.Cells.Item(36, 5).Value = "UP"
$ws.Cells.Item(36, 6).Value = 0.48
$ws.Cells.Item(36, 7).Value = ""
$ws.Cells.Item(36, 8).Value = "OPEN"
$ws.Cells.Item(36, 9).Value = 0
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 11).Value = 100
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 13).Value = 0
$ws.Cells.Item(36, 14).Value = 0.6
$ws.Cells.Item(36, 15).Value = "Normal spread capture: 19600 bps"
$ws.Cells.Item(36, 16).Value = ""
$ws.Cells.Item(36, 17).Value = 0

# ---- "MarketMaking" sheet: currently A1:Q2 -> new trade lands on row 3 ----
$ws2 = $wb.Worksheets.Item("MarketMaking")

$ws2.Cells.Item(3, 1).Value = 35

$ws2.Cells.Item(3, 2).NumberFormat = "@"
$ws2.Cells.Item(3, 2).Value = "2026-02-17"
$ws2.Cells.Item(3, 2).ClearFormats()

$ws2.Cells.Item(3, 3).Value = "20:17:26"
$ws2.Cells.Item(3, 4).Value = "MarketMaking"
$ws2.Cells.Item(3, 5).Value = "UP"
$ws2.Cells.Item(3, 6).Value = 0.48
$ws2.Cells.Item(3, 7).Value = ""
$ws2.Cells.Item(3, 8).Value = "OPEN"
$ws2.Cells.Item(3, 9).Value = 0
$ws2.Cells.Item(3, 10).Value = 0
$ws2.Cells.Item(3, 11).Value = 100
$ws2.Cells.Item(3, 12).Value = 0
$ws2.Cells.Item(3, 13).Value = 0
$ws2.Cells.Item(3, 14).Value = 0.6
$ws2.Cells.Item(3, 15).Value = "Normal spread capture: 19600 bps"
$ws2.Cells.Item(3, 16).Value = ""
$ws2.Cells.Item(3, 17).Value = 0
